# Insert a new data row at row 32 (pushing the existing rows 32-141 down to
# 33-142) and populate it with a new "Fruta / Plátano" price record, matching
# the weekly update described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 32..141 down to 33..142, carrying formatting along.
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new record.
$ws.Range("A32").Value = 1
$ws.Range("B32").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C32").Value = "Arica y Parinacota"
$ws.Range("D32").Value = 44481
$ws.Range("E32").Value = 15
$ws.Range("F32").Value = "Fruta"
$ws.Range("G32").Value = 100108
$ws.Range("H32").Value = "Tropicales y subtropicales"
$ws.Range("I32").Value = 100108006
$ws.Range("J32").Value = "Plátano"
$ws.Range("K32").Value = "Sin especificar"
$ws.Range("L32").Value = "Pintón"
$ws.Range("M32").Value = 120
$ws.Range("N32").Value = 20000
$ws.Range("O32").Value = 21000
$ws.Range("P32").Value = 20500
$ws.Range("Q32").Value = "$/caja 20 kilos"
$ws.Range("R32").Value = "Bolivia"
$ws.Range("S32").Value = 1025
$ws.Range("T32").Value = 20
